$d = $word.ActiveDocument

function Replace-ExactText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# NOTE: "910÷3=303, 1" is both a source (original cell text) and a destination
# (the replacement text for "153÷3=51, 0"). To avoid the newly-written
# "910÷3=303, 1" being matched by a later/overlapping replace, perform the
# replacement that consumes the original "910÷3=303, 1" text BEFORE writing
# a new "910÷3=303, 1" elsewhere.
Replace-ExactText "910÷3=303, 1" "397÷6=66, 1"

Replace-ExactText "344÷4=86, 0" "411÷5=82, 1"
Replace-ExactText "153÷3=51, 0" "910÷3=303, 1"
Replace-ExactText "858÷7=122, 4" "328÷9=36, 4"
Replace-ExactText "222÷8=27, 6" "476÷8=59, 4"
Replace-ExactText "759÷2=379, 1" "898÷3=299, 1"
Replace-ExactText "162÷6=27, 0" "727÷6=121, 1"
Replace-ExactText "225÷9=25, 0" "587÷2=293, 1"
Replace-ExactText "534÷8=66, 6" "290÷8=36, 2"
Replace-ExactText "675÷7=96, 3" "607÷3=202, 1"
Replace-ExactText "173÷2=86, 1" "136÷9=15, 1"
Replace-ExactText "525÷7=75, 0" "218÷2=109, 0"
Replace-ExactText "525÷8=65, 5" "222÷9=24, 6"
Replace-ExactText "998÷7=142, 4" "340÷7=48, 4"
Replace-ExactText "562÷2=281, 0" "808÷3=269, 1"
Replace-ExactText "568÷2=284, 0" "277÷5=55, 2"
Replace-ExactText "343÷6=57, 1" "194÷7=27, 5"
Replace-ExactText "251÷6=41, 5" "266÷5=53, 1"
Replace-ExactText "298÷5=59, 3" "868÷5=173, 3"
Replace-ExactText "736÷9=81, 7" "323÷5=64, 3"
Replace-ExactText "945÷7=135, 0" "998÷8=124, 6"
Replace-ExactText "495÷8=61, 7" "282÷9=31, 3"
Replace-ExactText "748÷8=93, 4" "719÷5=143, 4"
Replace-ExactText "702÷9=78, 0" "335÷5=67, 0"
Replace-ExactText "784÷6=130, 4" "164÷3=54, 2"
